# Trim excess trailing whitespace from cell values and fix several
# data-entry typos/format inconsistencies (e.g. "SSP DF" -> "SSP/DF",
# missing CPF/RG digits, renamed models, etc.) across the vehicle
# registry sheet, per the associated data-cleanup commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Placa "
$ws.Range("B1").Value = "Modelo "
$ws.Range("C1").Value = "Marca "
$ws.Range("D1").Value = "Cor "
$ws.Range("E1").Value = "Nome "
$ws.Range("F1").Value = "CPF "
$ws.Range("G1").Value = "RG "
# Row 2
$ws.Range("B2").Value = "Civic "
$ws.Range("C2").Value = "Honda "
$ws.Range("D2").Value = "Prata "
$ws.Range("E2").Value = "Gilveano Cota "
$ws.Range("F2").Value = "020.856.860-36 "
$ws.Range("G2").Value = "1092477941 SSP/RS "
# Row 3
$ws.Range("B3").Value = "HB20S "
$ws.Range("C3").Value = "Hyundai "
$ws.Range("D3").Value = "Prata "
$ws.Range("F3").Value = "016.016.291-25 "
$ws.Range("G3").Value = "2617795 SSP/DF "
# Row 4
$ws.Range("A4").Value = "RBU2C69 "
$ws.Range("B4").Value = "Onix "
$ws.Range("D4").Value = "Branco "
$ws.Range("E4").Value = "Rafaela Kuhn Valandro "
$ws.Range("F4").Value = "016.718.370-29 "
$ws.Range("G4").Value = "1077638623 SJS/FI RS "
# Row 5
$ws.Range("B5").Value = "Omega "
$ws.Range("D5").Value = "Preto "
$ws.Range("E5").Value = "Lucio José Assis da Silva "
$ws.Range("F5").Value = "054.172.355-39 "
$ws.Range("G5").Value = "33323500 SSP/SE "
# Row 6
$ws.Range("B6").Value = "Uno "
$ws.Range("C6").Value = "Fiat "
$ws.Range("D6").Value = "Cinza Escuro "
$ws.Range("E6").Value = "Thiago Santos Barros "
$ws.Range("F6").Value = "036.055.841-01 "
$ws.Range("G6").Value = "Ilegível "
# Row 7
$ws.Range("A7").Value = "RUK2E39 "
$ws.Range("B7").Value = "Mobi "
$ws.Range("C7").Value = "Fiat "
$ws.Range("D7").Value = "Cinza Escuro "
$ws.Range("E7").Value = "Solange Maria da Silva "
$ws.Range("F7").Value = "668.200.573-34 "
$ws.Range("G7").Value = "1545861 SSP/PI "
# Row 8
$ws.Range("A8").Value = ""
$ws.Range("C8").Value = "Jeep "
$ws.Range("D8").Value = "Branco "
$ws.Range("E8").Value = "Rafael Ramos Peres "
$ws.Range("F8").Value = "091.058.389-75 "
$ws.Range("G8").Value = "12430705-8 SESP/PR "
# Row 9
$ws.Range("A9").Value = ""
$ws.Range("B9").Value = "Creta "
$ws.Range("C9").Value = "Hyundai "
$ws.Range("D9").Value = "Preto "
$ws.Range("F9").Value = "020.494.546-19 "
$ws.Range("G9").Value = "MG18795190 SSP/MG "
# Row 10
$ws.Range("B10").Value = "Escort "
$ws.Range("C10").Value = "Ford "
$ws.Range("D10").Value = "Preto "
$ws.Range("E10").Value = "Otmar dos Reis Saffier "
$ws.Range("F10").Value = "726.944.757-87 "
# Row 11
$ws.Range("A11").Value = "PBF5129 "
$ws.Range("B11").Value = "Argo "
$ws.Range("C11").Value = "Fiat "
$ws.Range("D11").Value = "Prata "
$ws.Range("F11").Value = "011.728.465-37 "
$ws.Range("G11").Value = "1441805 SSP/SE "
# Row 12
$ws.Range("A12").Value = "PAH1E54 "
$ws.Range("B12").Value = "Ka "
$ws.Range("C12").Value = "Ford "
$ws.Range("D12").Value = "Branco "
$ws.Range("F12").Value = "466.264.681-20 "
$ws.Range("G12").Value = "462772 SSP/MS "
# Row 13
$ws.Range("A13").Value = "PBI-7069 "
$ws.Range("B13").Value = "Ranger "
$ws.Range("C13").Value = "Ford "
$ws.Range("D13").Value = "Cinza Escuro "
$ws.Range("E13").Value = "Fábio Chagas Theophilo "
$ws.Range("F13").Value = "Ilegível "
$ws.Range("G13").Value = "Ilegível "
# Row 14
$ws.Range("A14").Value = "REQ5D29 "
$ws.Range("B14").Value = "Tiguan "
$ws.Range("D14").Value = "Branco "
$ws.Range("E14").Value = "Benedito Sousa Alves "
$ws.Range("F14").Value = "054.739.813-6 "
$ws.Range("G14").Value = "036102752008-6 GE/UF/MA "
# Row 15
$ws.Range("A15").Value = "KDS9X69 "
$ws.Range("B15").Value = "Mobi "
$ws.Range("C15").Value = "Fiat "
$ws.Range("D15").Value = "Cinza Escuro "
$ws.Range("F15").Value = "219.634.508-33 "
$ws.Range("G15").Value = "3235320 SSP/SP "
# Row 16
$ws.Range("B16").Value = "Q3 "
$ws.Range("C16").Value = "Audi "
$ws.Range("D16").Value = "Preto "
$ws.Range("E16").Value = "Mauricio Carvalho Barros "
$ws.Range("F16").Value = "075.253.107-77 "
$ws.Range("G16").Value = "1081919821/PRJ "
# Row 17
$ws.Range("A17").Value = "OLV3G67 "
$ws.Range("B17").Value = "Palio "
$ws.Range("C17").Value = "Fiat "
$ws.Range("D17").Value = "Prata "
$ws.Range("F17").Value = "077.227.653-60 "
$ws.Range("G17").Value = "4138680 SSP/PI "
# Row 18
$ws.Range("B18").Value = "Fit "
$ws.Range("C18").Value = "Honda "
$ws.Range("D18").Value = "Prata "
$ws.Range("E18").Value = "José Aparecido Ribeiro "
$ws.Range("F18").Value = "061.969.978-74 "
$ws.Range("G18").Value = "16569506 SSP/SP "
# Row 19
$ws.Range("A19").Value = "SSU3186 "
$ws.Range("B19").Value = "Kwid "
$ws.Range("C19").Value = "Renault "
$ws.Range("D19").Value = "Branco "
$ws.Range("E19").Value = "Pedro Marini Lopes "
$ws.Range("F19").Value = "781.649.338-91 "
$ws.Range("G19").Value = "76632218 SSP/SP "
# Row 20
$ws.Range("D20").Value = "Preto "
$ws.Range("F20").Value = "714.085.849-72 "
$ws.Range("G20").Value = "22862249 SSP/SP "
# Row 21
$ws.Range("B21").Value = "Onix "
$ws.Range("D21").Value = "Branco "
$ws.Range("F21").Value = "161.639.748-30 "
$ws.Range("G21").Value = "10255175 SSP/SP "
# Row 22
$ws.Range("A22").Value = "GKJ3B09 "
$ws.Range("B22").Value = "Chevette "
$ws.Range("D22").Value = "Cinza "
$ws.Range("F22").Value = "766.879.712-04 "
$ws.Range("G22").Value = "415617 SSP/AC "
# Row 23
$ws.Range("A23").Value = "SGY1A99 "
$ws.Range("B23").Value = "Commander "
$ws.Range("C23").Value = "Jeep "
$ws.Range("D23").Value = "Prata "
$ws.Range("E23").Value = "Alan Machado Correa "
$ws.Range("F23").Value = "125.974.457-47 "
$ws.Range("G23").Value = "20817030DETRANRJ "
# Row 24
$ws.Range("A24").Value = "PRQ8F00 "
$ws.Range("E24").Value = "Luiz Fernando Ramos "
$ws.Range("F24").Value = "118.453.886-74 "
$ws.Range("G24").Value = "MG18846936 SSP "
